$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

# Update the recursive-filter coefficients
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1

# Update the data grid (rows 2-9, columns D:I)
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0

$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 0

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0

$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 0

$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0

$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0

# Update the active-cell selection
$ws.Range("B4").Select()
